$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two numeric values in row 2
$ws.Range("C2").Value = 2.263000011444092
$ws.Range("D2").Value = 20.0

# Clear out the "unused" helper columns C/D for rows 3-9 (formatting cleanup)
$ws.Range("C3:D9").ClearContents()

# Move the active selection to D2
$ws.Range("D2").Select() | Out-Null
